$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.400.13'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.642.27'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '591.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -1.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.642.21'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.62'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.44'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.121.17'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.320.04'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000144'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.640.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.37'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.74'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.20%  '
$ws.Range('E25').Value = '  +4.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.165'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.85'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '532.60'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +16.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0804'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '174.44'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.88'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +9.59%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('E40').Value = '  +6.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '172.03'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.40%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.28'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.72'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.11'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0557'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.629'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0960'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0238'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.58%  '
$ws.Range('E51').Value = '  -0.70%  '
